# Insert a new data row at row 275 (pushes the existing rows 275-380 down
# to 276-381), then populate the newly inserted row with its own data.
# This matches the commit's net effect: a new weekly price observation
# was added for "Macroferia Regional de Talca - Zanahoria", and every
# subsequent row shifted down by one (dimension grows from R380 to R381).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(275).EntireRow.Insert()

$ws.Range("A275").Value = 5
$ws.Range("B275").Value = "Macroferia Regional de Talca"
$ws.Range("C275").Value = "Maule"
$ws.Range("D275").Value = 44784
$ws.Range("E275").Value = 7
$ws.Range("F275").Value = 100114013
$ws.Range("G275").Value = "Zanahoria"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 500
$ws.Range("K275").Value = 9000
$ws.Range("L275").Value = 9000
$ws.Range("M275").Value = 9000
$ws.Range("N275").Value = "`$/saco 20 kilos"
$ws.Range("O275").Value = "Región de Ñuble"
$ws.Range("P275").Value = 450
$ws.Range("Q275").Value = 20
$ws.Range("R275").Value = "Hortaliza"
